$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the old "Inflammatory-Mac" rows (former rows 4 and 5), then
# rewrite all data rows (2-5) with the refreshed TPM-derived values.
# Simplest robust approach: delete rows 4 and 5 entirely (shifts former
# rows 6-7 up into rows 4-5), then overwrite every data cell with the
# final values from the new dataset.

$ws.Rows.Item(4).Delete() | Out-Null
$ws.Rows.Item(4).Delete() | Out-Null

# Row 2: FAPs -> Ereg -> Erbb4 -> FAPs
$ws.Cells.Item(2,1).Value = "FAPs"
$ws.Cells.Item(2,2).Value = "Ereg"
$ws.Cells.Item(2,3).Value = "Erbb4"
$ws.Cells.Item(2,4).Value = "FAPs"
$ws.Cells.Item(2,5).Value = 1
$ws.Cells.Item(2,6).Value = 0.3333333333333333
$ws.Cells.Item(2,7).Value = 0.1313356666666667
$ws.Cells.Item(2,8).Value = 0.394007
$ws.Cells.Item(2,9).Value = 0.9277061342889635
$ws.Cells.Item(2,10).Value = 0.9277061342889635
$ws.Cells.Item(2,11).Value = 1
$ws.Cells.Item(2,12).Value = 0.3333333333333333
$ws.Cells.Item(2,13).Value = 0.010299
$ws.Cells.Item(2,14).Value = 0.030897
$ws.Cells.Item(2,15).Value = 0.9810128591839975
$ws.Cells.Item(2,16).Value = 0.9810128591839974
$ws.Cells.Item(2,17).Value = 0.001352626031
$ws.Cells.Item(2,18).Value = 0.012173634279
$ws.Cells.Item(2,19).Value = 0.9100916472813496
$ws.Cells.Item(2,20).Value = 0.9100916472813495

# Row 3: FAPs -> Ereg -> Erbb4 -> MuSCs
$ws.Cells.Item(3,1).Value = "FAPs"
$ws.Cells.Item(3,2).Value = "Ereg"
$ws.Cells.Item(3,3).Value = "Erbb4"
$ws.Cells.Item(3,4).Value = "MuSCs"
$ws.Cells.Item(3,5).Value = 1
$ws.Cells.Item(3,6).Value = 0.3333333333333333
$ws.Cells.Item(3,7).Value = 0.1313356666666667
$ws.Cells.Item(3,8).Value = 0.394007
$ws.Cells.Item(3,9).Value = 0.9277061342889635
$ws.Cells.Item(3,10).Value = 0.9277061342889635
$ws.Cells.Item(3,11).Value = 1
$ws.Cells.Item(3,12).Value = 0.3333333333333333
$ws.Cells.Item(3,13).Value = 0.0001993333333333333
$ws.Cells.Item(3,14).Value = 0.000598
$ws.Cells.Item(3,15).Value = 0.01898714081600254
$ws.Cells.Item(3,16).Value = 0.01898714081600254
$ws.Cells.Item(3,17).Value = 0.00002617957622222222
$ws.Cells.Item(3,18).Value = 0.000235616186
$ws.Cells.Item(3,19).Value = 0.01761448700761391
$ws.Cells.Item(3,20).Value = 0.01761448700761391

# Row 4: MuSCs -> Ereg -> Erbb4 -> FAPs
$ws.Cells.Item(4,1).Value = "MuSCs"
$ws.Cells.Item(4,2).Value = "Ereg"
$ws.Cells.Item(4,3).Value = "Erbb4"
$ws.Cells.Item(4,4).Value = "FAPs"
$ws.Cells.Item(4,5).Value = 1
$ws.Cells.Item(4,6).Value = 0.3333333333333333
$ws.Cells.Item(4,7).Value = 0.01023466666666667
$ws.Cells.Item(4,8).Value = 0.030704
$ws.Cells.Item(4,9).Value = 0.07229386571103645
$ws.Cells.Item(4,10).Value = 0.07229386571103644
$ws.Cells.Item(4,11).Value = 1
$ws.Cells.Item(4,12).Value = 0.3333333333333333
$ws.Cells.Item(4,13).Value = 0.010299
$ws.Cells.Item(4,14).Value = 0.030897
$ws.Cells.Item(4,15).Value = 0.9810128591839975
$ws.Cells.Item(4,16).Value = 0.9810128591839974
$ws.Cells.Item(4,17).Value = 0.000105406832
$ws.Cells.Item(4,18).Value = 0.000948661488
$ws.Cells.Item(4,19).Value = 0.07092121190264783
$ws.Cells.Item(4,20).Value = 0.0709212119026478

# Row 5: MuSCs -> Ereg -> Erbb4 -> MuSCs
$ws.Cells.Item(5,1).Value = "MuSCs"
$ws.Cells.Item(5,2).Value = "Ereg"
$ws.Cells.Item(5,3).Value = "Erbb4"
$ws.Cells.Item(5,4).Value = "MuSCs"
$ws.Cells.Item(5,5).Value = 1
$ws.Cells.Item(5,6).Value = 0.3333333333333333
$ws.Cells.Item(5,7).Value = 0.01023466666666667
$ws.Cells.Item(5,8).Value = 0.030704
$ws.Cells.Item(5,9).Value = 0.07229386571103645
$ws.Cells.Item(5,10).Value = 0.07229386571103644
$ws.Cells.Item(5,11).Value = 1
$ws.Cells.Item(5,12).Value = 0.3333333333333333
$ws.Cells.Item(5,13).Value = 0.0001993333333333333
$ws.Cells.Item(5,14).Value = 0.000598
$ws.Cells.Item(5,15).Value = 0.01898714081600254
$ws.Cells.Item(5,16).Value = 0.01898714081600254
$ws.Cells.Item(5,17).Value = 0.000002040110222222222
$ws.Cells.Item(5,18).Value = 0.000018360992
$ws.Cells.Item(5,19).Value = 0.001372653808388627
$ws.Cells.Item(5,20).Value = 0.001372653808388626
